# T4 and Stylesheet Work
# Improved handling of <lb> elements so that the line number is displayed
# in all contexts -> add a new glyph entry ("g50" / "apostrophe") to the
# "Glyphs" lookup sheet, right after the existing "g49" / "k abbreviation"
# row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Glyphs")

# Append the new row (row 51) below the existing data (which ends at row 50).
$ws.Range("A51").Value = "g50"
$ws.Range("B51").Value = "apostrophe"

# Leave the selection where the author left it when they saved the file.
$ws.Range("E52").Select()
